# Weekly CompStat report refresh: new crime data collected.
# Updates the report header (volume/issue number, week-covering date range)
# and rewrites the crime-complaint statistics table (rows 14-30) with the
# newly collected figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1. Header text tweaks (surgical, run-preserving replacements)
# ---------------------------------------------------------------------

# "Volume 30   Number  6"  ->  "Volume 30   Number  7"
$volCell = $ws.Range("A8")
$volText = $volCell.Value()
$numIdx = $volText.IndexOf("6") + 1
$volCell.Characters($numIdx, 1).Text = "7"

# "Report Covering the Week  2/6/2023  Through  2/12/2023"
#   -> "Report Covering the Week  2/13/2023  Through  2/19/2023"
$weekCell = $ws.Range("C9")
$weekText = $weekCell.Value()
$d1Idx = $weekText.IndexOf("2/6/2023") + 1
$weekCell.Characters($d1Idx, 8).Text = "2/13/2023"
$weekText = $weekCell.Value()
$d2Idx = $weekText.IndexOf("2/12/2023") + 1
$weekCell.Characters($d2Idx, 9).Text = "2/19/2023"

# ---------------------------------------------------------------------
# 2. Crime statistics table (columns C..N, rows 14-30)
# ---------------------------------------------------------------------
# Each row value is either a number, or one of the two placeholder
# strings used by the report when a ratio is not meaningful:
#   "0"      -> used in the (normally numeric) count columns
#   "***.*"  -> used in the (normally numeric) % change columns
#
# Column order used below: C, D, E, F, G, H, I, J, K, L, M, N

$rowData = @{
  14 = @("0", "0", "***.*", "0", 1, -100, "0", 2, -100, -100, -100, -100)
  15 = @("0", 2, -100, "0", 3, -100, 1, 4, -75, -66.666666666666, -50, -91.666666666666)
  16 = @(5, 6, -16.666666666666, 19, 12, 58.333333333333, 30, 20, 50, 100, -28.571428571428, -82.352941176470)
  17 = @(3, 4, -25, 13, 20, -35, 33, 38, -13.157894736842, 6.451612903225, -28.260869565217, -68.571428571428)
  18 = @(2, 10, -80, 12, 25, -52, 25, 31, -19.354838709677, 13.636363636363, 4.166666666666, -69.879518072289)
  19 = @(5, 8, -37.5, 26, 31, -16.129032258064, 53, 49, 8.163265306122, 35.897435897435, 82.758620689655, 20.454545454545)
  20 = @(2, "0", "***.*", 9, 15, -40, 13, 24, -45.833333333333, 0, 0, -84.883720930232)
  21 = @(17, 30, -43.333333333333, 79, 107, -26.168224299065, 155, 168, -7.738095238095, 25, -2.515723270440, -69.428007889546)
  22 = @("0", 1, -100, "0", 2, -100, 1, 4, -75, -66.666666666666, 0, "***.*")
  23 = @(1, 1, 0, 5, 7, -28.571428571428, 12, 11, 9.090909090909, 9.090909090909, 0, "***.*")
  24 = @(23, 12, 91.666666666666, 73, 51, 43.137254901960, 118, 86, 37.209302325581, 93.442622950819, 40.476190476190, "***.*")
  25 = @(11, 9, 22.222222222222, 33, 31, 6.451612903225, 56, 53, 5.660377358490, 64.705882352941, -42.268041237113, "***.*")
  26 = @("0", 2, -100, 1, 3, -66.666666666666, 3, 5, -40, -40, "***.*", "***.*")
  27 = @("0", "0", "***.*", 1, 2, -50, 1, 5, -80, -87.5, "***.*", "***.*")
  28 = @("0", "0", "***.*", "0", 1, -100, 2, 3, -33.333333333333, 100, -77.777777777777, -92)
  29 = @("0", "0", "***.*", "0", 1, -100, 2, 3, -33.333333333333, 100, -77.777777777777, -91.304347826087)
  30 = @("0", "0", "***.*", "0", "0", "***.*", "0", "0", "***.*", "***.*", "***.*", "***.*")
}

$cols = @("C", "D", "E", "F", "G", "H", "I", "J", "K", "L", "M", "N")
# Columns that hold a percentage-style figure (and so use the "***.*"
# placeholder / percent number format) versus a plain count (which uses
# the "0" placeholder / integer number format).
$percentCols = @("E", "H", "K", "L", "M", "N")

# Style template cells, taken from row 14 which already exhibits every
# flavor of formatting used in the data rows (text placeholder, integer,
# percentage). Row 21 (TOTAL) keeps its own bold styles untouched since
# none of its cells change type there.
$templateText = $ws.Range("C14")
$templateInt = $ws.Range("D14")
$templatePct = $ws.Range("H14")

foreach ($r in 14..30) {
    $values = $rowData[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $col = $cols[$i]
        $value = $values[$i]
        $cell = $ws.Range("$col$r")

        if ($value -is [string]) {
            # Text placeholder ("0" or "***.*"). Force text number format
            # first so Excel does not silently coerce "0" back into a
            # numeric zero, then restore the proper (General) style.
            $cell.NumberFormat = "@"
            $cell.Value = $value
            $templateText.Copy()
            $cell.PasteSpecial($xlPasteFormats)
        } else {
            $isPercent = $percentCols -contains $col
            if ($isPercent) {
                $templatePct.Copy()
            } else {
                $templateInt.Copy()
            }
            $cell.PasteSpecial($xlPasteFormats)
            $cell.Value = $value
        }
    }
}

$excel.CutCopyMode = $false
